# nominate_ambassadors.xlsx edit script
# Commit: "added field as example to verify nominator knows nominee.
#          made sure various elements are required."
#
# Logical change on the "survey" worksheet (sheet1):
#  - The nominee_name calculate row (formerly row 19, right after nominee_key)
#    is moved down, below the update_phone / nominee_name_new / nominee_phone_new
#    question block (which moves up to directly follow nominee_key).
#  - A brand-new "nominee_age" integer question is inserted right after
#    nominee_phone_new, with a constraint requiring an age between 18 and 110.
#  - nominee_name_new, nominee_phone_new, and nominee_age are all marked
#    "required" (column K = "yes").
#  - The nominee_phone calculate formula gets its spacing normalized
#    (${nominee_new}=1 -> ${nominee_new} = 1, etc.) as a byproduct of
#    re-entry, matching the committed text exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# ---------------------------------------------------------------------
# 1) Restructure rows: pull the "nominee_name" calculate row out from
#    right after nominee_key, then re-create it after the phone block,
#    leaving room for the new "nominee_age" question in between.
# ---------------------------------------------------------------------

# Row 19 currently holds: calculate / nominee_name / ... (the one we're moving)
$ws.Rows.Item(19).Delete()

# After the delete, row order is:
#   17 nominee_new, 18 nominee_key, 19 update_phone, 20 nominee_name_new,
#   21 nominee_phone_new, 22 nominee_phone(calc), 23 nominee_label(calc), ...
# Insert two fresh blank rows at 22/23: one for the new nominee_age question,
# one to re-host the nominee_name calculate row we removed above.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(23).Insert()

# Fix up formatting on the two newly-inserted blank rows, since a bare
# Insert() doesn't fully inherit borders from neighboring rows. Copy
# formatting (not values) from a same-shaped neighbor row.
$ws.Range("A21:N21").Copy() | Out-Null
$ws.Range("A22:N22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A24:N24").Copy() | Out-Null
$ws.Range("A23:N23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Set cell content for the re-shuffled block (rows 19-24).
# ---------------------------------------------------------------------

# Row 19: select_one yesno / update_phone (unchanged content, new position)
$ws.Range("A19").Value2 = 'select_one yesno'
$ws.Range("B19").Value2 = 'update_phone'
$ws.Range("C19").Value2 = 'Would you like to provide ${nominee_name}''s phone number?'
$ws.Range("I19").Value2 = '${nominee_new} = 0 and pulldata(''nominees'', ''nominee_phone'', ''nominee_key'', ${nominee_key_new}) = -1'
$ws.Range("K19").Value2 = 'yes'

# Row 20: text / nominee_name_new (unchanged content, new position; now required)
$ws.Range("A20").Value2 = 'text'
$ws.Range("B20").Value2 = 'nominee_name_new'
$ws.Range("C20").Value2 = 'Enter the nominee''s name.'
$ws.Range("I20").Value2 = '${nominee_new} = 1'
$ws.Range("K20").Value2 = 'yes'

# Row 21: integer / nominee_phone_new (unchanged content, new position; now required)
$ws.Range("A21").Value2 = 'integer'
$ws.Range("B21").Value2 = 'nominee_phone_new'
$ws.Range("C21").Value2 = 'Enter ${nominee_name}''s phone number, or enter -1 if the phone number is unknown or cannot be provided.'
$ws.Range("G21").Value2 = '. = -1 or (. > 99 and coalesce(pulldata(''nominees'', ''nominee_phone'', ''nominee_phone'', .), 1) = 1)'
$ws.Range("H21").Value2 = 'That number is invalid or corresponds to an existing nominee.'
$ws.Range("I21").Value2 = '${nominee_new} = 1 or ${update_phone} = 1'
$ws.Range("K21").Value2 = 'yes'

# Row 22: NEW integer / nominee_age question, required, age sanity constraint.
$ws.Range("A22").Value2 = 'integer'
$ws.Range("B22").Value2 = 'nominee_age'
$ws.Range("C22").Value2 = 'Enter ${nominee_name}''s approximate age in years.'
$ws.Range("G22").Value2 = '. >= 18 and . <= 110'
$ws.Range("K22").Value2 = 'yes'

# Row 23: calculate / nominee_name (re-created here, after the phone block)
$ws.Range("A23").Value2 = 'calculate'
$ws.Range("B23").Value2 = 'nominee_name'
$ws.Range("N23").Value2 = 'if(${nominee_new} = 1, ${nominee_name_new}, pulldata(''nominees'', ''nominee_name'', ''nominee_key'', ${nominee_key}))'

# Row 24: calculate / nominee_phone (content re-entered with normalized spacing)
$ws.Range("A24").Value2 = 'calculate'
$ws.Range("B24").Value2 = 'nominee_phone'
$ws.Range("N24").Value2 = 'if(${nominee_new} = 1 or ${update_phone} = 1, ${nominee_phone_new}, pulldata(''nominees'', ''nominee_phone'', ''nominee_key'', ${nominee_key}))'

# ---------------------------------------------------------------------
# 3) Row heights, to match the committed layout.
# ---------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 51
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(21).RowHeight = 68
$ws.Rows.Item(22).RowHeight = 34
$ws.Rows.Item(31).RowHeight = 34

# ---------------------------------------------------------------------
# 4) Selection / frozen-pane view state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("K24").Select()

Write-Output "done"
